$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> @{ D = newPrice; E = newVolume }  (only cells that actually changed)
$updates = @{
    2  = @{ D = "26.910.05"; E = "  +0.02%  " }
    3  = @{ D = "1.815.20";  E = "  +0.33%  " }
    5  = @{ D = "309.00";    E = "  -0.47%  " }
    6  = @{             E = "  +0.03%  " }
    7  = @{             E = "  +0.87%  " }
    8  = @{ D = "0.3658";    E = "  -1.21%  " }
    9  = @{ D = "0.07360";   E = "  -0.20%  " }
    10 = @{ D = "0.8685";    E = "  -0.57%  " }
    11 = @{ D = "20.22";     E = "  -1.19%  " }
    12 = @{ D = "1.801.72";  E = "  -0.13%  " }
    13 = @{ D = "5.375";     E = "  +0.35%  " }
    14 = @{ D = "0.07095";   E = "  +0.66%  " }
    15 = @{ D = "6.497";     E = "  -0.51%  " }
    16 = @{ D = "91.15";     E = "  -0.74%  " }
    17 = @{             E = "  +0.10%  " }
    18 = @{ D = "0.000008679"; E = "  -0.11%  " }
    20 = @{ D = "14.63";     E = "  -0.59%  " }
    21 = @{ D = "26.930.07"; E = "  +0.10%  " }
    22 = @{ D = "5.294";     E = "  -0.55%  " }
    23 = @{             E = "  -0.90%  " }
    24 = @{ D = "2.045.73";  E = "  +1.08%  " }
    25 = @{ D = "1.896";     E = "  +0.00%  " }
    26 = @{ D = "150.96";    E = "  -0.15%  " }
    27 = @{ D = "18.39";     E = "  +0.20%  " }
    28 = @{ D = "2.135";     E = "  -0.28%  " }
    29 = @{ D = "5.263";     E = "  -0.68%  " }
    30 = @{ D = "115.83";    E = "  -0.05%  " }
    31 = @{ D = "0.08902";   E = "  +0.01%  " }
    32 = @{ D = "0.7566";    E = "  +0.50%  " }
    33 = @{ D = "1.164";     E = "  +0.79%  " }
    34 = @{ D = "4.479";     E = "  +0.79%  " }
    35 = @{ D = "2.902";     E = "  -0.29%  " }
    36 = @{             E = "  +0.06%  " }
    37 = @{             E = "  -1.01%  " }
    38 = @{ D = "0.05276";   E = "  +0.83%  " }
    39 = @{ D = "0.01945";   E = "  -1.30%  " }
    40 = @{ D = "2.964";     E = "  +1.13%  " }
    41 = @{ D = "7.176";     E = "  +0.05%  " }
    42 = @{ D = "0.5280";    E = "  -0.18%  " }
    43 = @{ D = "2.337";     E = "  -3.58%  " }
    44 = @{             E = "  -0.50%  " }
    45 = @{ D = "8.424";     E = "  -0.92%  " }
    46 = @{ D = "0.4842";    E = "  -2.53%  " }
    47 = @{             E = "  +0.59%  " }
    48 = @{             E = "  +0.06%  " }
    49 = @{ D = "103.16";    E = "  -0.75%  " }
    50 = @{             E = "  -0.69%  " }
    51 = @{ D = "0.06292";   E = "  -0.05%  " }
}

foreach ($row in $updates.Keys) {
    $cellChanges = $updates[$row]

    if ($cellChanges.ContainsKey("D")) {
        # The Price column holds plain text (e.g. "26.910.05", "309.00",
        # "0.3658") that Excel's COM Value setter would otherwise silently
        # coerce into a number (dropping trailing zeros / changing
        # precision). Force the cell to Text first, write the literal
        # string, then restore the Normal style so no stray formatting is
        # left behind on the cell.
        $dCell = $ws.Range("D$row")
        $dCell.NumberFormat = "@"
        $dCell.Value = $cellChanges["D"]
        $dCell.Style = "Normal"
    }

    if ($cellChanges.ContainsKey("E")) {
        # Volume(1h) values (e.g. "  +0.02%  ") are never numeric-looking,
        # so a plain Value assignment keeps them as text with the padding
        # spaces intact.
        $ws.Range("E$row").Value = $cellChanges["E"]
    }
}
